# burlington.xlsx - choropleth steps update / muni-code cleanup
#
# This reproduces the column-width / hidden-state changes and the
# selection change captured in the source diff for sheet
# "BURLINGTON COUNTY 2022". A handful of columns that used to hold
# helper/hidden data (O:U, Y:Z) are unhidden and resized now that they
# feed the choropleth steps, column K (11) is widened, and the active
# selection is moved to P1 after scrolling the sheet to the newly
# relevant rows near the bottom of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width / visibility changes -------------------------------
# NOTE: the COM ColumnWidth setter in this host quantizes to the
# nearest 1/6 of a character (it stores width = ColumnWidth + 5/6),
# so the values below are chosen to land as close as possible on the
# target "number of characters" widths from the authored workbook.

# Column K (11): widen from ~43.43 to 51 characters
$ws.Columns.Item(11).ColumnWidth = 50.166666666666664

# Column O (15): unhide + resize (was hidden helper column)
$ws.Columns.Item(15).ColumnWidth = 10.333333333333332
$ws.Columns.Item(15).Hidden = $false

# Column P (16): unhide + resize
$ws.Columns.Item(16).ColumnWidth = 12.5
$ws.Columns.Item(16).Hidden = $false

# Column Q (17): unhide + resize
$ws.Columns.Item(17).ColumnWidth = 12.666666666666666
$ws.Columns.Item(17).Hidden = $false

# Column R (18): unhide + resize
$ws.Columns.Item(18).ColumnWidth = 19.166666666666668
$ws.Columns.Item(18).Hidden = $false

# Column S (19): unhide + resize
$ws.Columns.Item(19).ColumnWidth = 17.666666666666668
$ws.Columns.Item(19).Hidden = $false

# Column T (20): unhide + resize
$ws.Columns.Item(20).ColumnWidth = 13.166666666666666
$ws.Columns.Item(20).Hidden = $false

# Column U (21): unhide + resize
$ws.Columns.Item(21).ColumnWidth = 17.666666666666668
$ws.Columns.Item(21).Hidden = $false

# Column Y (25): resize (already visible)
$ws.Columns.Item(25).ColumnWidth = 24.666666666666668

# Column Z (26): unhide + resize (was hidden, width 0)
$ws.Columns.Item(26).ColumnWidth = 9.666666666666666
$ws.Columns.Item(26).Hidden = $false

# --- Selection ---------------------------------------------------------
# Scroll near the bottom of the table and leave the active selection on
# P1, matching the saved view state in the workbook.
$excel.ActiveWindow.ScrollRow = 247
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("P1").Select()
